# Update gh-pages to output generated at 456a3b4
# Applies:
#  - Sheet "展览" (index 1): insert two new rows (new row 28 "上海·多厨狂喜动漫展1.0"
#    and new row 37 "上海·NW新界动漫游戏展"), shifting following rows down, plus
#    various F-column ("想去人数") bumps.
#  - Sheet "演出" (index 2): F-column bumps only.
#  - Sheet "本地生活" (index 3): F-column bumps only.
#  - Sheet "全部类型" (index 4): F-column bumps only.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# --- insert new row 28: 上海·多厨狂喜动漫展1.0 -----------------------------
$ws1.Rows.Item(28).Insert()
# copy formatting (bold / border / alignment) from the row above so the new
# row's A-column index cell matches the sheet's existing style
$ws1.Cells.Item(27,1).Copy() | Out-Null
$ws1.Cells.Item(28,1).PasteSpecial(-4122) | Out-Null

$ws1.Cells.Item(28,1).Value = 27
$ws1.Cells.Item(28,2).NumberFormat = "@"
$ws1.Cells.Item(28,2).Value = "2024-05-03"
$ws1.Cells.Item(28,3).Value = "上海·多厨狂喜动漫展1.0"
$ws1.Cells.Item(28,4).Value = "澳门路168号 月星家居（澳门路）"
$ws1.Cells.Item(28,5).Value = "2024.05.03 10:00-05.04 16:00"
$ws1.Cells.Item(28,6).Value = 0
$ws1.Cells.Item(28,7).Value = 49
$ws1.Cells.Item(28,8).Value = "https://show.bilibili.com/platform/detail.html?id=83932"
$ws1.Cells.Item(28,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/721HW21G1712565123950.jpeg"

# --- insert new row 37: 上海·NW新界动漫游戏展 ------------------------------
$ws1.Rows.Item(37).Insert()
$ws1.Cells.Item(36,1).Copy() | Out-Null
$ws1.Cells.Item(37,1).PasteSpecial(-4122) | Out-Null

$ws1.Cells.Item(37,1).Value = 36
$ws1.Cells.Item(37,2).NumberFormat = "@"
$ws1.Cells.Item(37,2).Value = "2024-06-08"
$ws1.Cells.Item(37,3).Value = "上海·NW新界动漫游戏展"
$ws1.Cells.Item(37,4).Value = "长寿路街道澳门路168号 月星家居"
$ws1.Cells.Item(37,5).Value = "2024.06.08 10:00-06.10 16:00"
$ws1.Cells.Item(37,6).Value = 0
$ws1.Cells.Item(37,7).Value = 49
$ws1.Cells.Item(37,8).Value = "https://show.bilibili.com/platform/detail.html?id=83934"
$ws1.Cells.Item(37,9).Value = "//i0.hdslb.com/bfs/openplatform/202404/4eW55lpD1712517166770.jpeg"

# --- fix up the manual running-index column A for every row after the first
#     insertion point; Excel does NOT recompute literal (non-formula) values
#     in shifted rows, so rows 29-38 keep their pre-shift A value unless we
#     rewrite them (A holds row-number-minus-one throughout this sheet).
for ($r = 29; $r -le 38; $r++) {
    $ws1.Cells.Item($r, 1).Value = $r - 1
}

# --- F-column ("想去人数") bumps on existing rows (post-insert row numbers) -
$sheet1Updates = @{
    3 = 785; 4 = 549; 5 = 2272; 6 = 1353; 8 = 807; 11 = 2969; 14 = 1089;
    17 = 222; 19 = 1035; 20 = 1035; 21 = 133; 22 = 13; 23 = 162; 25 = 195; 26 = 634;
    29 = 823; 30 = 49; 32 = 1022; 33 = 5009; 34 = 473; 35 = 223
}
foreach ($r in $sheet1Updates.Keys) {
    $ws1.Cells.Item($r, 6).Value = $sheet1Updates[$r]
}

# ---------------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$sheet2Updates = @{
    19 = 4; 23 = 48; 24 = 373; 26 = 646; 28 = 2; 34 = 416; 35 = 416; 41 = 748; 42 = 40
}
foreach ($r in $sheet2Updates.Keys) {
    $ws2.Cells.Item($r, 6).Value = $sheet2Updates[$r]
}

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$sheet3Updates = @{ 4 = 634; 5 = 410; 6 = 394 }
foreach ($r in $sheet3Updates.Keys) {
    $ws3.Cells.Item($r, 6).Value = $sheet3Updates[$r]
}

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$sheet4Updates = @{
    4 = 410; 5 = 785; 7 = 549; 9 = 2272; 10 = 1353; 12 = 807; 16 = 2969; 19 = 1089;
    22 = 394; 26 = 222; 27 = 1035; 28 = 1035; 29 = 133; 30 = 4; 31 = 13; 32 = 162;
    33 = 195; 35 = 634; 37 = 373; 38 = 646; 39 = 823; 40 = 49; 41 = 1022; 42 = 5009;
    44 = 473; 46 = 416; 47 = 223; 51 = 40
}
foreach ($r in $sheet4Updates.Keys) {
    $ws4.Cells.Item($r, 6).Value = $sheet4Updates[$r]
}

Write-Output "edit complete"
